$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '48.244.84'
$ws.Range('E2').Value = '  +2.03%  '
$ws.Range('D3').Value = '2.511.42'
$ws.Range('E3').Value = '  +0.84%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '109.38'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.99%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '320.67'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.07%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.530'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.74%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.999'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.03%  '
$ws.Range('E9').Value = '  +2.42%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '40.00'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +3.06%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '20.18'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +10.06%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0819'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.16%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.125'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.92%  '
$ws.Range('E14').Value = '  +1.17%  '
$ws.Range('D15').Value = '2.903.13'
$ws.Range('E15').Value = '  +0.81%  '
$ws.Range('D16').Value = '2.508.24'
$ws.Range('E16').Value = '  +0.78%  '
$ws.Range('E17').Value = '  +0.37%  '
$ws.Range('D18').Value = '48.096.52'
$ws.Range('E18').Value = '  +1.90%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.21'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.01%  '
$ws.Range('E20').Value = '  +0.04%  '
$ws.Range('D21').Value = '0.0₃0945'
$ws.Range('E21').Value = '  +1.28%  '
$ws.Range('E22').Value = '  +2.10%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '72.21'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.76%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '274.73'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +12.10%  '
$ws.Range('E25').Value = '  -0.32%  '
$ws.Range('E26').Value = '  -0.02%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '25.98'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.08%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.39'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +4.91%  '
$ws.Range('E29').Value = '  +0.64%  '
$ws.Range('B30').Value = 'InjectiveProtocol'
$ws.Range('C30').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '35.54'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +2.66%  '
$ws.Range('B31').Value = 'Kaspa'
$ws.Range('C31').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.140'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +2.25%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '49.38'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.51%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '19.36'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -6.80%  '
$ws.Range('E34').Value = '  +0.20%  '
$ws.Range('E35').Value = '  -0.01%  '
$ws.Range('E37').Value = '  +0.17%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.65'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.22%  '
$ws.Range('E39').Value = '  +1.20%  '
$ws.Range('E40').Value = '  +1.02%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '122.15'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +3.63%  '
$ws.Range('E42').Value = '  -0.90%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '21.85'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -6.36%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0305'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +3.31%  '
$ws.Range('D45').Value = '2.031.80'
$ws.Range('E45').Value = '  +2.08%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.13'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +3.31%  '
$ws.Range('E47').Value = '  +5.28%  '
$ws.Range('E48').Value = '  -1.53%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '9.04'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.96%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '5.19'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +2.04%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '79.86'
$ws.Range('D51').Style = 'Normal'
